$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("paper")

$ws.Cells.Item(2, 3).Value = "(2018, 2.9934292653876886)"
$ws.Cells.Item(2, 5).Value = 222.0304411
$ws.Cells.Item(2, 6).Value = -0.1082029
$ws.Cells.Item(2, 7).Value = -395.1055604
$ws.Cells.Item(2, 8).Value = 0.5077121
$ws.Cells.Item(2, 9).Value = -0.0001537
$ws.Cells.Item(2, 10).Value = 222.0304411
$ws.Cells.Item(2, 11).Value = -0.1082029
$ws.Cells.Item(2, 12).Value = -395.1055604
$ws.Cells.Item(2, 13).Value = 0.5077121
$ws.Cells.Item(2, 14).Value = -0.0001537

$ws.Cells.Item(3, 3).Value = "(2018, 4.335035967360962)"
$ws.Cells.Item(3, 5).Value = 247.3717352
$ws.Cells.Item(3, 6).Value = -0.1206782
$ws.Cells.Item(3, 7).Value = 3769.942885
$ws.Cells.Item(3, 8).Value = -3.6362795
$ws.Cells.Item(3, 9).Value = 0.0008771
$ws.Cells.Item(3, 10).Value = 247.3717352
$ws.Cells.Item(3, 11).Value = -0.1206782
$ws.Cells.Item(3, 12).Value = 3769.942885
$ws.Cells.Item(3, 13).Value = -3.6362795
$ws.Cells.Item(3, 14).Value = 0.0008771

$ws.Cells.Item(4, 3).Value = "(2018, 7.726415385162614)"
$ws.Cells.Item(4, 5).Value = -67.2460648
$ws.Cells.Item(4, 6).Value = 0.037236
$ws.Cells.Item(4, 7).Value = 6487.9641441
$ws.Cells.Item(4, 8).Value = -6.505004
$ws.Cells.Item(4, 9).Value = 0.0016323
$ws.Cells.Item(4, 10).Value = -67.2460648
$ws.Cells.Item(4, 11).Value = 0.037236
$ws.Cells.Item(4, 12).Value = 6487.9641441
$ws.Cells.Item(4, 13).Value = -6.505004
$ws.Cells.Item(4, 14).Value = 0.0016323

$ws.Cells.Item(5, 2).Value = $true
$ws.Cells.Item(5, 3).Value = "(2018, 2.3831305034058814)"
$ws.Cells.Item(5, 5).Value = 298.6321488
$ws.Cells.Item(5, 6).Value = -0.1466218
$ws.Cells.Item(5, 7).Value = -31874.6613446
$ws.Cells.Item(5, 8).Value = 31.9630132
$ws.Cells.Item(5, 9).Value = -0.0080114
$ws.Cells.Item(5, 10).Value = 298.6321488
$ws.Cells.Item(5, 11).Value = -0.1466218
$ws.Cells.Item(5, 12).Value = -31874.6613446
$ws.Cells.Item(5, 13).Value = 31.9630132
$ws.Cells.Item(5, 14).Value = -0.0080114

$ws.Cells.Item(6, 3).Value = "(2018, 3.3125209564987204)"
$ws.Cells.Item(6, 5).Value = 98.8448157
$ws.Cells.Item(6, 6).Value = -0.0472717
$ws.Cells.Item(6, 7).Value = -8705.3324614
$ws.Cells.Item(6, 8).Value = 8.7394855
$ws.Cells.Item(6, 9).Value = -0.0021923
$ws.Cells.Item(6, 10).Value = 98.8448157
$ws.Cells.Item(6, 11).Value = -0.0472717
$ws.Cells.Item(6, 12).Value = -8705.3324614
$ws.Cells.Item(6, 13).Value = 8.7394855
$ws.Cells.Item(6, 14).Value = -0.0021923

$ws.Cells.Item(7, 3).Value = "(2018, 15.3594)"

$ws.Cells.Item(8, 3).Value = "(2018, 1.4207488018932446)"

$ws.Cells.Item(9, 3).Value = "(2018, 3.2388130004710316)"

$ws.Cells.Item(10, 3).Value = "(2018, 3.6335709562563583)"

$ws.Cells.Item(11, 3).Value = "(2018, 3.350555669589186)"
$ws.Cells.Item(11, 5).Value = 153.4365599
$ws.Cells.Item(11, 6).Value = -0.0749659
$ws.Cells.Item(11, 7).Value = 49589.8098032
$ws.Cells.Item(11, 8).Value = -49.4135237
$ws.Cells.Item(11, 9).Value = 0.01231
$ws.Cells.Item(11, 10).Value = 153.4365599
$ws.Cells.Item(11, 11).Value = -0.0749659
$ws.Cells.Item(11, 12).Value = 49589.8098032
$ws.Cells.Item(11, 13).Value = -49.4135237
$ws.Cells.Item(11, 14).Value = 0.01231

$ws.Cells.Item(12, 3).Value = "(2018, 3.36198828541822)"
$ws.Cells.Item(12, 5).Value = 136.2618055
$ws.Cells.Item(12, 6).Value = -0.0659043
$ws.Cells.Item(12, 7).Value = 6839.9587035
$ws.Cells.Item(12, 8).Value = -6.7563372
$ws.Cells.Item(12, 9).Value = 0.0016693
$ws.Cells.Item(12, 10).Value = 136.2618055
$ws.Cells.Item(12, 11).Value = -0.0659043
$ws.Cells.Item(12, 12).Value = 6839.9587035
$ws.Cells.Item(12, 13).Value = -6.7563372
$ws.Cells.Item(12, 14).Value = 0.0016693

$ws.Cells.Item(13, 3).Value = "(2018, 4.939284601531085)"
$ws.Cells.Item(13, 5).Value = 1369.6938774
$ws.Cells.Item(13, 6).Value = -0.6795178
$ws.Cells.Item(13, 7).Value = 185115.8737055
$ws.Cells.Item(13, 8).Value = -184.0621352
$ws.Cells.Item(13, 9).Value = 0.0457541
$ws.Cells.Item(13, 10).Value = 1369.6938774
$ws.Cells.Item(13, 11).Value = -0.6795178
$ws.Cells.Item(13, 12).Value = 185115.8737055
$ws.Cells.Item(13, 13).Value = -184.0621352
$ws.Cells.Item(13, 14).Value = 0.0457541

$ws.Cells.Item(14, 3).Value = "(2014, 26.07954329210276)"
$ws.Cells.Item(14, 5).Value = -67.4980105
$ws.Cells.Item(14, 6).Value = 0.0425624
$ws.Cells.Item(14, 7).Value = -966966.9661009
$ws.Cells.Item(14, 8).Value = 964.0557577
$ws.Cells.Item(14, 9).Value = -0.2402825
$ws.Cells.Item(14, 10).Value = -67.4980105
$ws.Cells.Item(14, 11).Value = 0.0425624
$ws.Cells.Item(14, 12).Value = -966966.9661009
$ws.Cells.Item(14, 13).Value = 964.0557577
$ws.Cells.Item(14, 14).Value = -0.2402825

$ws.Cells.Item(15, 3).Value = "(2018, 3.777473176726379)"
$ws.Cells.Item(15, 5).Value = 9.8864548
$ws.Cells.Item(15, 6).Value = -0.0029807
$ws.Cells.Item(15, 7).Value = 23826.8255832
$ws.Cells.Item(15, 8).Value = -23.7745934
$ws.Cells.Item(15, 9).Value = 0.0059315
$ws.Cells.Item(15, 10).Value = 9.8864548
$ws.Cells.Item(15, 11).Value = -0.0029807
$ws.Cells.Item(15, 12).Value = 23826.8255832
$ws.Cells.Item(15, 13).Value = -23.7745934
$ws.Cells.Item(15, 14).Value = 0.0059315

$ws.Cells.Item(16, 3).Value = "(2018, 2.714738255033557)"
$ws.Cells.Item(16, 5).Value = 69.3158592
$ws.Cells.Item(16, 6).Value = -0.0327117
$ws.Cells.Item(16, 7).Value = -18488.6066069
$ws.Cells.Item(16, 8).Value = 18.4884919
$ws.Cells.Item(16, 9).Value = -0.0046211
$ws.Cells.Item(16, 10).Value = 69.3158592
$ws.Cells.Item(16, 11).Value = -0.0327117
$ws.Cells.Item(16, 12).Value = -18488.6066069
$ws.Cells.Item(16, 13).Value = 18.4884919
$ws.Cells.Item(16, 14).Value = -0.0046211

$ws.Cells.Item(17, 2).Value = $true
$ws.Cells.Item(17, 3).Value = "(2018, 3.4382535691607106)"
$ws.Cells.Item(17, 5).Value = 78.3873866
$ws.Cells.Item(17, 6).Value = -0.0372328
$ws.Cells.Item(17, 7).Value = 9585.8469996
$ws.Cells.Item(17, 8).Value = -9.5258808
$ws.Cells.Item(17, 9).Value = 0.0023674
$ws.Cells.Item(17, 10).Value = 78.3873866
$ws.Cells.Item(17, 11).Value = -0.0372328
$ws.Cells.Item(17, 12).Value = 9585.8469996
$ws.Cells.Item(17, 13).Value = -9.5258808
$ws.Cells.Item(17, 14).Value = 0.0023674

$ws.Cells.Item(18, 3).Value = "(2018, 3.446236534782528)"
$ws.Cells.Item(18, 5).Value = 262.5846277
$ws.Cells.Item(18, 6).Value = -0.1287047
$ws.Cells.Item(18, 7).Value = 14048.1923624
$ws.Cells.Item(18, 8).Value = -13.8870361
$ws.Cells.Item(18, 9).Value = 0.0034327
$ws.Cells.Item(18, 10).Value = 262.5846277
$ws.Cells.Item(18, 11).Value = -0.1287047
$ws.Cells.Item(18, 12).Value = 14048.1923624
$ws.Cells.Item(18, 13).Value = -13.8870361
$ws.Cells.Item(18, 14).Value = 0.0034327

$ws.Cells.Item(19, 3).Value = "(2018, 5.332684821125188)"
$ws.Cells.Item(19, 5).Value = 108.2446001
$ws.Cells.Item(19, 6).Value = -0.0511042
$ws.Cells.Item(19, 7).Value = 1482.9506685
$ws.Cells.Item(19, 8).Value = -1.4230902
$ws.Cells.Item(19, 9).Value = 0.0003423
$ws.Cells.Item(19, 10).Value = 108.2446001
$ws.Cells.Item(19, 11).Value = -0.0511042
$ws.Cells.Item(19, 12).Value = 1482.9506685
$ws.Cells.Item(19, 13).Value = -1.4230902
$ws.Cells.Item(19, 14).Value = 0.0003423

$ws.Cells.Item(20, 2).Value = $true
$ws.Cells.Item(20, 3).Value = "(2018, 4.662113749083906)"
$ws.Cells.Item(20, 5).Value = 181.6121292
$ws.Cells.Item(20, 6).Value = -0.0872691
$ws.Cells.Item(20, 7).Value = 65865.135861
$ws.Cells.Item(20, 8).Value = -65.6408306
$ws.Cells.Item(20, 9).Value = 0.0163557
$ws.Cells.Item(20, 10).Value = 181.6121292
$ws.Cells.Item(20, 11).Value = -0.0872691
$ws.Cells.Item(20, 12).Value = 65865.135861
$ws.Cells.Item(20, 13).Value = -65.6408306
$ws.Cells.Item(20, 14).Value = 0.0163557

$ws.Cells.Item(21, 3).Value = "(2018, 3.4097014163583004)"
$ws.Cells.Item(21, 5).Value = 121.5174218
$ws.Cells.Item(21, 6).Value = -0.0587189
$ws.Cells.Item(21, 7).Value = 4868.6540898
$ws.Cells.Item(21, 8).Value = -4.7964628
$ws.Cells.Item(21, 9).Value = 0.0011821
$ws.Cells.Item(21, 10).Value = 121.5174218
$ws.Cells.Item(21, 11).Value = -0.0587189
$ws.Cells.Item(21, 12).Value = 4868.6540898
$ws.Cells.Item(21, 13).Value = -4.7964628
$ws.Cells.Item(21, 14).Value = 0.0011821

$ws.Cells.Item(22, 3).Value = "(2018, 3.3929890617925316)"
$ws.Cells.Item(22, 5).Value = -172.6628968
$ws.Cells.Item(22, 6).Value = 0.0880174
$ws.Cells.Item(22, 7).Value = -68960.3574879
$ws.Cells.Item(22, 8).Value = 68.7311852
$ws.Cells.Item(22, 9).Value = -0.0171244
$ws.Cells.Item(22, 10).Value = -172.6628968
$ws.Cells.Item(22, 11).Value = 0.0880174
$ws.Cells.Item(22, 12).Value = -68960.3574879
$ws.Cells.Item(22, 13).Value = 68.7311852
$ws.Cells.Item(22, 14).Value = -0.0171244

$ws.Cells.Item(23, 3).Value = "(2018, 6.242981790591806)"
$ws.Cells.Item(23, 5).Value = 112.8326837
$ws.Cells.Item(23, 6).Value = -0.0528196
$ws.Cells.Item(23, 7).Value = -10128.6417037
$ws.Cells.Item(23, 8).Value = 10.1683909
$ws.Cells.Item(23, 9).Value = -0.0025502
$ws.Cells.Item(23, 10).Value = 112.8326837
$ws.Cells.Item(23, 11).Value = -0.0528196
$ws.Cells.Item(23, 12).Value = -10128.6417037
$ws.Cells.Item(23, 13).Value = 10.1683909
$ws.Cells.Item(23, 14).Value = -0.0025502

$ws.Cells.Item(24, 3).Value = "(2018, 7.1105413667291195)"

$ws.Cells.Item(25, 2).Value = $true
$ws.Cells.Item(25, 3).Value = "(2018, 9.598252696456086)"
$ws.Cells.Item(25, 5).Value = -726.9843291
$ws.Cells.Item(25, 6).Value = 0.3672272
$ws.Cells.Item(25, 7).Value = 46962.9572053
$ws.Cells.Item(25, 8).Value = -47.2283544
$ws.Cells.Item(25, 9).Value = 0.0118751
$ws.Cells.Item(25, 10).Value = -726.9843291
$ws.Cells.Item(25, 11).Value = 0.3672272
$ws.Cells.Item(25, 12).Value = 46962.9572053
$ws.Cells.Item(25, 13).Value = -47.2283544
$ws.Cells.Item(25, 14).Value = 0.0118751

$ws.Cells.Item(26, 3).Value = "(2018, 12.274977954144621)"
$ws.Cells.Item(26, 5).Value = 86.3493511
$ws.Cells.Item(26, 6).Value = -0.0371144
$ws.Cells.Item(26, 7).Value = 54055.6774833
$ws.Cells.Item(26, 8).Value = -53.8996582
$ws.Cells.Item(26, 9).Value = 0.0134388
$ws.Cells.Item(26, 10).Value = 86.3493511
$ws.Cells.Item(26, 11).Value = -0.0371144
$ws.Cells.Item(26, 12).Value = 54055.6774833
$ws.Cells.Item(26, 13).Value = -53.8996582
$ws.Cells.Item(26, 14).Value = 0.0134388

$ws.Cells.Item(29, 3).Value = "(2018, 1.7222916666666668)"
$ws.Cells.Item(29, 5).Value = 309.0283042
$ws.Cells.Item(29, 6).Value = -0.1528262
$ws.Cells.Item(29, 7).Value = 15662.2193942
$ws.Cells.Item(29, 8).Value = -15.4737594
$ws.Cells.Item(29, 9).Value = 0.0038221
$ws.Cells.Item(29, 10).Value = 309.0283042
$ws.Cells.Item(29, 11).Value = -0.1528262
$ws.Cells.Item(29, 12).Value = 15662.2193942
$ws.Cells.Item(29, 13).Value = -15.4737594
$ws.Cells.Item(29, 14).Value = 0.0038221

$ws.Cells.Item(31, 3).Value = "(2018, 2.6970269299820466)"
$ws.Cells.Item(31, 5).Value = -404.1827282
$ws.Cells.Item(31, 6).Value = 0.2021817
$ws.Cells.Item(31, 7).Value = -56721.467369
$ws.Cells.Item(31, 8).Value = 56.2847924
$ws.Cells.Item(31, 9).Value = -0.0139621
$ws.Cells.Item(31, 10).Value = -404.1827282
$ws.Cells.Item(31, 11).Value = 0.2021817
$ws.Cells.Item(31, 12).Value = -56721.467369
$ws.Cells.Item(31, 13).Value = 56.2847924
$ws.Cells.Item(31, 14).Value = -0.0139621

$ws.Cells.Item(33, 3).Value = "(2018, 12.247575757575758)"

$ws.Cells.Item(34, 3).Value = "(2018, 4.059476707083599)"
$ws.Cells.Item(34, 5).Value = 529.8635718
$ws.Cells.Item(34, 6).Value = -0.2612381
$ws.Cells.Item(34, 7).Value = 24665.4398856
$ws.Cells.Item(34, 8).Value = -24.3490594
$ws.Cells.Item(34, 9).Value = 0.0060099
$ws.Cells.Item(34, 10).Value = 529.8635718
$ws.Cells.Item(34, 11).Value = -0.2612381
$ws.Cells.Item(34, 12).Value = 24665.4398856
$ws.Cells.Item(34, 13).Value = -24.3490594
$ws.Cells.Item(34, 14).Value = 0.0060099

$ws.Cells.Item(35, 3).Value = "(2018, 22.62025811497849)"
$ws.Cells.Item(35, 5).Value = -1373.5975932
$ws.Cells.Item(35, 6).Value = 0.6910903
$ws.Cells.Item(35, 7).Value = 87722.6693273
$ws.Cells.Item(35, 8).Value = -88.2288896
$ws.Cells.Item(35, 9).Value = 0.0221856
$ws.Cells.Item(35, 10).Value = -1373.5975932
$ws.Cells.Item(35, 11).Value = 0.6910903
$ws.Cells.Item(35, 12).Value = 87722.6693273
$ws.Cells.Item(35, 13).Value = -88.2288896
$ws.Cells.Item(35, 14).Value = 0.0221856
